# Generate Report for Handoff
# Adds a new tracked file "b0589368-4664-441c-96e1-fc884252e46d.md" to the
# localization status report, inserting a row for it just above the
# always-last ".localization-config" row on all three worksheets
# (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newFile   = "b0589368-4664-441c-96e1-fc884252e46d.md"
$newHash   = "d8c12315abfa683eac90066995103cc8b736479e"
$zhFile    = "$newFile.$newHash.zh-cn.xlf"
$deFile    = "$newFile.$newHash.de-de.xlf"
$zhDate    = "2016-03-07 09:41:35"
$deDate    = "2016-03-07 09:41:47"
$epoch     = "0001-01-01 00:00:00"

$newFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$newHash/e2e/$newFile"
$configUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/2ce524d1b44aa53b3f19812b0b96a9a49c0a4352/.localization-config"
$zhXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$newHash/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhFile"
$deXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$newHash/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deFile"

# ---------------------------------------------------------------------
# Sheet 1: "Overview" -- 3 columns (File Name, zh-cn, de-de)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Rows.Item(9).Insert()

$ws1.Range("A9").Value = $newFile
$ws1.Range("B9").Value = "Ready for handoff"
$ws1.Range("C9").Value = "Ready for handoff"

$ws1.Range("A10").Value = ".localization-config"
$ws1.Range("B10").Value = "Not to be localized"
$ws1.Range("C10").Value = "Not to be localized"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce524d1b44aa53b3f19812b0b96a9a49c0a4352/e2e/fc27a2e6-d0cf-487c-9349-c5c89fa5783c.md", "", "", "fc27a2e6-d0cf-487c-9349-c5c89fa5783c.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/239d0879e57d1b1f7251a4fead4d24d0ac7c73ff/e2e/06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.md", "", "", "06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/66f5884563424c7539484faee238f819a2af811d/e2e/785ed964-5661-4e1a-aa6f-dc1102863f4f.md", "", "", "785ed964-5661-4e1a-aa6f-dc1102863f4f.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/239d0879e57d1b1f7251a4fead4d24d0ac7c73ff/e2e/dbb8b7c2-6efa-4117-90da-56923094cd06.md", "", "", "dbb8b7c2-6efa-4117-90da-56923094cd06.md")
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/f4055970680395dc3b7eecf6d9ddfc37e7948de6/e2e/fc368083-54a6-4157-a90d-0fc2c1a1fe53.md", "", "", "fc368083-54a6-4157-a90d-0fc2c1a1fe53.md")
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/9c78ddd6f28012dea28437b984c3a8b47508a7b3/e2e/05294cb7-e2f2-411b-a2bd-ca4347d00657.md", "", "", "05294cb7-e2f2-411b-a2bd-ca4347d00657.md")
$ws1.Hyperlinks.Add($ws1.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/20b0f243690f693b27e01c51db9dd3f6ff589aa5/e2e/ae39c42a-b758-451a-8e19-8fafa59f77a6.md", "", "", "ae39c42a-b758-451a-8e19-8fafa59f77a6.md")
$ws1.Hyperlinks.Add($ws1.Range("A9"), $newFileUrl, "", "", $newFile)
$ws1.Hyperlinks.Add($ws1.Range("A10"), $configUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn" -- 9 columns
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Rows.Item(9).Insert()

$ws2.Range("A9").Value = $newFile
$ws2.Range("B9").Value = "Ready for handoff"
$ws2.Range("C9").Value = $zhFile
$ws2.Range("D9").Value = $zhDate
$ws2.Range("G9").Value = $epoch
$ws2.Range("H9").Value = "Include"

$ws2.Range("A10").Value = ".localization-config"
$ws2.Range("B10").Value = "Not to be localized"
$ws2.Range("D10").Value = $epoch
$ws2.Range("G10").Value = $epoch
$ws2.Range("H10").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce524d1b44aa53b3f19812b0b96a9a49c0a4352/e2e/fc27a2e6-d0cf-487c-9349-c5c89fa5783c.md", "", "", "fc27a2e6-d0cf-487c-9349-c5c89fa5783c.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8b648f6d7dad889edad9851689dc20ea272d5277/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/fc27a2e6-d0cf-487c-9349-c5c89fa5783c.6a2d1e28dc288a5beca5b4a27939d9814d36e661.zh-cn.xlf", "", "", "fc27a2e6-d0cf-487c-9349-c5c89fa5783c.6a2d1e28dc288a5beca5b4a27939d9814d36e661.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/545fab8b1a24f44259799bc2b1d18a723a1c0159/e2e/fc27a2e6-d0cf-487c-9349-c5c89fa5783c.md", "", "", "fc27a2e6-d0cf-487c-9349-c5c89fa5783c.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5b6657060c0175554b674b7bbb707a81ce0f4335/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/fc27a2e6-d0cf-487c-9349-c5c89fa5783c.6a2d1e28dc288a5beca5b4a27939d9814d36e661.zh-cn.xlf", "", "", "fc27a2e6-d0cf-487c-9349-c5c89fa5783c.6a2d1e28dc288a5beca5b4a27939d9814d36e661.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/239d0879e57d1b1f7251a4fead4d24d0ac7c73ff/e2e/06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.md", "", "", "06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f566921d63ecb759a7405dc81c59ef598af8399/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.0e70ba733ba7709f7ac6be56e782aba445bc8037.zh-cn.xlf", "", "", "06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.0e70ba733ba7709f7ac6be56e782aba445bc8037.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/66f5884563424c7539484faee238f819a2af811d/e2e/785ed964-5661-4e1a-aa6f-dc1102863f4f.md", "", "", "785ed964-5661-4e1a-aa6f-dc1102863f4f.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71016509f163159583c1a2b9df87ae678a593d7b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/785ed964-5661-4e1a-aa6f-dc1102863f4f.073894905ffb5f94af9424482a49a95b8aec6e26.zh-cn.xlf", "", "", "785ed964-5661-4e1a-aa6f-dc1102863f4f.073894905ffb5f94af9424482a49a95b8aec6e26.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b41813478354d1bd19fc41722bbd5cb60e35cb9b/e2e/785ed964-5661-4e1a-aa6f-dc1102863f4f.md", "", "", "785ed964-5661-4e1a-aa6f-dc1102863f4f.md")
$ws2.Hyperlinks.Add($ws2.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/97a2d8b5a6e63de81ac9eebfb6e387d22bf35d47/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/785ed964-5661-4e1a-aa6f-dc1102863f4f.073894905ffb5f94af9424482a49a95b8aec6e26.zh-cn.xlf", "", "", "785ed964-5661-4e1a-aa6f-dc1102863f4f.073894905ffb5f94af9424482a49a95b8aec6e26.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/239d0879e57d1b1f7251a4fead4d24d0ac7c73ff/e2e/dbb8b7c2-6efa-4117-90da-56923094cd06.md", "", "", "dbb8b7c2-6efa-4117-90da-56923094cd06.md")
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f566921d63ecb759a7405dc81c59ef598af8399/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/dbb8b7c2-6efa-4117-90da-56923094cd06.9101628fb3d38b6c45007db9c130d59806c44cb6.zh-cn.xlf", "", "", "dbb8b7c2-6efa-4117-90da-56923094cd06.9101628fb3d38b6c45007db9c130d59806c44cb6.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/f4055970680395dc3b7eecf6d9ddfc37e7948de6/e2e/fc368083-54a6-4157-a90d-0fc2c1a1fe53.md", "", "", "fc368083-54a6-4157-a90d-0fc2c1a1fe53.md")
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c21722051eaa0afb7f3f87cf3650a9a99b61a757/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/fc368083-54a6-4157-a90d-0fc2c1a1fe53.17ffdc14645d6678fc3687045001f2fb09a8a6dc.zh-cn.xlf", "", "", "fc368083-54a6-4157-a90d-0fc2c1a1fe53.17ffdc14645d6678fc3687045001f2fb09a8a6dc.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/9c78ddd6f28012dea28437b984c3a8b47508a7b3/e2e/05294cb7-e2f2-411b-a2bd-ca4347d00657.md", "", "", "05294cb7-e2f2-411b-a2bd-ca4347d00657.md")
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/28aeefee9c0399befa5c7cf511ee7493db326010/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/05294cb7-e2f2-411b-a2bd-ca4347d00657.68f81ef12b8a67246ab4a8ce9182b2002eeb404d.zh-cn.xlf", "", "", "05294cb7-e2f2-411b-a2bd-ca4347d00657.68f81ef12b8a67246ab4a8ce9182b2002eeb404d.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/20b0f243690f693b27e01c51db9dd3f6ff589aa5/e2e/ae39c42a-b758-451a-8e19-8fafa59f77a6.md", "", "", "ae39c42a-b758-451a-8e19-8fafa59f77a6.md")
$ws2.Hyperlinks.Add($ws2.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2fa338be068c82e49a3bbbfa178738ca20827c6d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ae39c42a-b758-451a-8e19-8fafa59f77a6.8ff77a13ed90721b8fd70e62deaf3ce1a1491f1c.zh-cn.xlf", "", "", "ae39c42a-b758-451a-8e19-8fafa59f77a6.8ff77a13ed90721b8fd70e62deaf3ce1a1491f1c.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A9"), $newFileUrl, "", "", $newFile)
$ws2.Hyperlinks.Add($ws2.Range("C9"), $zhXlfUrl, "", "", $zhFile)
$ws2.Hyperlinks.Add($ws2.Range("A10"), $configUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 3: "de-de" -- 9 columns
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Rows.Item(9).Insert()

$ws3.Range("A9").Value = $newFile
$ws3.Range("B9").Value = "Ready for handoff"
$ws3.Range("C9").Value = $deFile
$ws3.Range("D9").Value = $deDate
$ws3.Range("G9").Value = $epoch
$ws3.Range("H9").Value = "Include"

$ws3.Range("A10").Value = ".localization-config"
$ws3.Range("B10").Value = "Not to be localized"
$ws3.Range("D10").Value = $epoch
$ws3.Range("G10").Value = $epoch
$ws3.Range("H10").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2ce524d1b44aa53b3f19812b0b96a9a49c0a4352/e2e/fc27a2e6-d0cf-487c-9349-c5c89fa5783c.md", "", "", "fc27a2e6-d0cf-487c-9349-c5c89fa5783c.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5accd3cd07b7496b091295acbd1dd284980cf30/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/fc27a2e6-d0cf-487c-9349-c5c89fa5783c.6a2d1e28dc288a5beca5b4a27939d9814d36e661.de-de.xlf", "", "", "fc27a2e6-d0cf-487c-9349-c5c89fa5783c.6a2d1e28dc288a5beca5b4a27939d9814d36e661.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fb84409d22a55a0c161366047144e3afaf8f5d24/e2e/fc27a2e6-d0cf-487c-9349-c5c89fa5783c.md", "", "", "fc27a2e6-d0cf-487c-9349-c5c89fa5783c.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9fbbe2141e9af373991f7aa57d24f678bebce904/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/fc27a2e6-d0cf-487c-9349-c5c89fa5783c.6a2d1e28dc288a5beca5b4a27939d9814d36e661.de-de.xlf", "", "", "fc27a2e6-d0cf-487c-9349-c5c89fa5783c.6a2d1e28dc288a5beca5b4a27939d9814d36e661.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/239d0879e57d1b1f7251a4fead4d24d0ac7c73ff/e2e/06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.md", "", "", "06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8b2e0e8513621865dd16d434a18b4bcf509d4fbc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.0e70ba733ba7709f7ac6be56e782aba445bc8037.de-de.xlf", "", "", "06f5abdb-2d3d-43a0-af9d-7fa2fc81f184.0e70ba733ba7709f7ac6be56e782aba445bc8037.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/66f5884563424c7539484faee238f819a2af811d/e2e/785ed964-5661-4e1a-aa6f-dc1102863f4f.md", "", "", "785ed964-5661-4e1a-aa6f-dc1102863f4f.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b6abb08865933b3f9fd38c5d0430a6e5933961e2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/785ed964-5661-4e1a-aa6f-dc1102863f4f.073894905ffb5f94af9424482a49a95b8aec6e26.de-de.xlf", "", "", "785ed964-5661-4e1a-aa6f-dc1102863f4f.073894905ffb5f94af9424482a49a95b8aec6e26.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/180f92ad29de7c5f62b0543eed9ceb0e33d620ef/e2e/785ed964-5661-4e1a-aa6f-dc1102863f4f.md", "", "", "785ed964-5661-4e1a-aa6f-dc1102863f4f.md")
$ws3.Hyperlinks.Add($ws3.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3f1d4e0683df7081fbabbcd9919639a85047cd82/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/785ed964-5661-4e1a-aa6f-dc1102863f4f.073894905ffb5f94af9424482a49a95b8aec6e26.de-de.xlf", "", "", "785ed964-5661-4e1a-aa6f-dc1102863f4f.073894905ffb5f94af9424482a49a95b8aec6e26.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/239d0879e57d1b1f7251a4fead4d24d0ac7c73ff/e2e/dbb8b7c2-6efa-4117-90da-56923094cd06.md", "", "", "dbb8b7c2-6efa-4117-90da-56923094cd06.md")
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8b2e0e8513621865dd16d434a18b4bcf509d4fbc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/dbb8b7c2-6efa-4117-90da-56923094cd06.9101628fb3d38b6c45007db9c130d59806c44cb6.de-de.xlf", "", "", "dbb8b7c2-6efa-4117-90da-56923094cd06.9101628fb3d38b6c45007db9c130d59806c44cb6.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/f4055970680395dc3b7eecf6d9ddfc37e7948de6/e2e/fc368083-54a6-4157-a90d-0fc2c1a1fe53.md", "", "", "fc368083-54a6-4157-a90d-0fc2c1a1fe53.md")
$ws3.Hyperlinks.Add($ws3.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b2db43b3d259b99f5ee21ca15eecddda9566fe4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/fc368083-54a6-4157-a90d-0fc2c1a1fe53.17ffdc14645d6678fc3687045001f2fb09a8a6dc.de-de.xlf", "", "", "fc368083-54a6-4157-a90d-0fc2c1a1fe53.17ffdc14645d6678fc3687045001f2fb09a8a6dc.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/9c78ddd6f28012dea28437b984c3a8b47508a7b3/e2e/05294cb7-e2f2-411b-a2bd-ca4347d00657.md", "", "", "05294cb7-e2f2-411b-a2bd-ca4347d00657.md")
$ws3.Hyperlinks.Add($ws3.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/413d8b59fa90258f9243ba0af727baca699eb31c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/05294cb7-e2f2-411b-a2bd-ca4347d00657.68f81ef12b8a67246ab4a8ce9182b2002eeb404d.de-de.xlf", "", "", "05294cb7-e2f2-411b-a2bd-ca4347d00657.68f81ef12b8a67246ab4a8ce9182b2002eeb404d.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/20b0f243690f693b27e01c51db9dd3f6ff589aa5/e2e/ae39c42a-b758-451a-8e19-8fafa59f77a6.md", "", "", "ae39c42a-b758-451a-8e19-8fafa59f77a6.md")
$ws3.Hyperlinks.Add($ws3.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/94ba7bf758b558e015b99d6d321a9942b457531b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ae39c42a-b758-451a-8e19-8fafa59f77a6.8ff77a13ed90721b8fd70e62deaf3ce1a1491f1c.de-de.xlf", "", "", "ae39c42a-b758-451a-8e19-8fafa59f77a6.8ff77a13ed90721b8fd70e62deaf3ce1a1491f1c.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A9"), $newFileUrl, "", "", $newFile)
$ws3.Hyperlinks.Add($ws3.Range("C9"), $deXlfUrl, "", "", $deFile)
$ws3.Hyperlinks.Add($ws3.Range("A10"), $configUrl, "", "", ".localization-config")

Write-Host "Done: inserted handoff row for $newFile on all sheets."
